$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting rows 3-10 down to rows 4-11
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 contents
$ws.Range("A3").Value = "app.fancy-slider.controls"
$ws.Range("B3").Style = "Good"
$ws.Range("C3").Style = "Good"
$ws.Range("D3").Value = "Trebuie integrat!"
$ws.Range("D3").Font.Color = 26012
$ws.Range("D3").Interior.Color = 10284031
$ws.Range("E3").Value = 'Ii comunica "Creierului" in ce directie vrea userul sa schimbe sliderul.'

# Update selection to match the diff
$ws.Range("E3").Select()
